$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The TPM recompute dropped the "ECs" target-cluster rows (old rows 2-4 as a target,
# plus rows 10-13 for the Resolving-Mac sending-cluster block); the table now has
# only 8 data rows (FAPs/MuSCs targets) instead of 12. Remove the trailing 4 rows
# first so row indices line up with the new data.
$ws.Rows("10:13").Delete()

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Adm"
$ws.Range("C2").Value = "Calcr"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 13.26223533333333
$ws.Range("H2").Value = 39.786706
$ws.Range("I2").Value = 0.391007655706778
$ws.Range("J2").Value = 0.3910076557067781
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.0005823333333333334
$ws.Range("N2").Value = 0.001747
$ws.Range("O2").Value = 0.00009625315715314126
$ws.Range("P2").Value = 0.00009625315715314125
$ws.Range("Q2").Value = 0.00772304170911111
$ws.Range("R2").Value = 0.069507375382
$ws.Range("S2").Value = 0.00003763572133282586
$ws.Range("T2").Value = 0.00003763572133282586

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Adm"
$ws.Range("C3").Value = "Calcr"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 13.26223533333333
$ws.Range("H3").Value = 39.786706
$ws.Range("I3").Value = 0.391007655706778
$ws.Range("J3").Value = 0.3910076557067781
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 6.049435666666667
$ws.Range("N3").Value = 18.148307
$ws.Range("O3").Value = 0.9999037468428469
$ws.Range("P3").Value = 0.9999037468428468
$ws.Range("Q3").Value = 80.22903944519354
$ws.Range("R3").Value = 722.0613550067419
$ws.Range("S3").Value = 0.3909700199854452
$ws.Range("T3").Value = 0.3909700199854452

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Adm"
$ws.Range("C4").Value = "Calcr"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 19.55844
$ws.Range("H4").Value = 58.67532
$ws.Range("I4").Value = 0.5766373150128344
$ws.Range("J4").Value = 0.5766373150128344
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.0005823333333333334
$ws.Range("N4").Value = 0.001747
$ws.Range("O4").Value = 0.00009625315715314126
$ws.Range("P4").Value = 0.00009625315715314125
$ws.Range("Q4").Value = 0.01138953156
$ws.Range("R4").Value = 0.10250578404
$ws.Range("S4").Value = 0.00005550316210229578
$ws.Range("T4").Value = 0.00005550316210229577

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Adm"
$ws.Range("C5").Value = "Calcr"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 19.55844
$ws.Range("H5").Value = 58.67532
$ws.Range("I5").Value = 0.5766373150128344
$ws.Range("J5").Value = 0.5766373150128344
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 6.049435666666667
$ws.Range("N5").Value = 18.148307
$ws.Range("O5").Value = 0.9999037468428469
$ws.Range("P5").Value = 0.9999037468428468
$ws.Range("Q5").Value = 118.31752452036
$ws.Range("R5").Value = 1064.85772068324
$ws.Range("S5").Value = 0.5765818118507322
$ws.Range("T5").Value = 0.5765818118507321

# Row 6
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Adm"
$ws.Range("C6").Value = "Calcr"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.8919193333333334
$ws.Range("H6").Value = 2.675758
$ws.Range("I6").Value = 0.02629626747232247
$ws.Range("J6").Value = 0.02629626747232247
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.0005823333333333334
$ws.Range("N6").Value = 0.001747
$ws.Range("O6").Value = 0.00009625315715314126
$ws.Range("P6").Value = 0.00009625315715314125
$ws.Range("Q6").Value = 0.0005193943584444445
$ws.Range("R6").Value = 0.004674549226
$ws.Range("S6").Value = 0.000002531098765554491
$ws.Range("T6").Value = 0.000002531098765554491

# Row 7
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Adm"
$ws.Range("C7").Value = "Calcr"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.8919193333333334
$ws.Range("H7").Value = 2.675758
$ws.Range("I7").Value = 0.02629626747232247
$ws.Range("J7").Value = 0.02629626747232247
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 6.049435666666667
$ws.Range("N7").Value = 18.148307
$ws.Range("O7").Value = 0.9999037468428469
$ws.Range("P7").Value = 0.9999037468428468
$ws.Range("Q7").Value = 5.395608626856222
$ws.Range("R7").Value = 48.560477641706
$ws.Range("S7").Value = 0.02629373637355691
$ws.Range("T7").Value = 0.02629373637355691

# Row 8
$ws.Range("A8").Value = "Resolving-Mac"
$ws.Range("B8").Value = "Adm"
$ws.Range("C8").Value = "Calcr"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 0.2055016666666667
$ws.Range("H8").Value = 0.6165050000000001
$ws.Range("I8").Value = 0.006058761808064916
$ws.Range("J8").Value = 0.006058761808064917
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.0005823333333333334
$ws.Range("N8").Value = 0.001747
$ws.Range("O8").Value = 0.00009625315715314126
$ws.Range("P8").Value = 0.00009625315715314125
$ws.Range("Q8").Value = 0.0001196704705555556
$ws.Range("R8").Value = 0.001077034235
$ws.Range("S8").Value = 0.0000005831749524651227
$ws.Range("T8").Value = 0.0000005831749524651227

# Row 9
$ws.Range("A9").Value = "Resolving-Mac"
$ws.Range("B9").Value = "Adm"
$ws.Range("C9").Value = "Calcr"
$ws.Range("D9").Value = "MuSCs"
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 0.2055016666666667
$ws.Range("H9").Value = 0.6165050000000001
$ws.Range("I9").Value = 0.006058761808064916
$ws.Range("J9").Value = 0.006058761808064917
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 6.049435666666667
$ws.Range("N9").Value = 18.148307
$ws.Range("O9").Value = 0.9999037468428469
$ws.Range("P9").Value = 0.9999037468428468
$ws.Range("Q9").Value = 1.243169111892778
$ws.Range("R9").Value = 11.188522007035
$ws.Range("S9").Value = 0.006058178633112451
$ws.Range("T9").Value = 0.006058178633112451
